# Add files via upload
# Populates the per-hotel/zone rate columns (B:G) for rows 2-30 on the only
# worksheet, and updates the active selection to C4 (single cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = [ordered]@{
    2  = @(12,5,10,10,3,3)
    3  = @(12,5,10,10,3,3)
    4  = @(16,5,10,10,3,3)
    5  = @(14,5,10,10,3,3)
    6  = @(12,5,10,10,3,3)
    7  = @(12,5,10,10,3,3)
    8  = @(12,5,10,10,3,3)
    9  = @(12,5,10,10,3,3)
    10 = @(12,5,10,10,3,3)
    11 = @(12,5,10,10,3,3)
    12 = @(4,5,10,10,3,3)
    13 = @(12,5,10,10,3,3)
    14 = @(12,5,10,10,3,3)
    15 = @(12,5,10,10,3,3)
    16 = @(12,5,10,10,3,3)
    17 = @(12,5,10,10,3,3)
    18 = @(3,3,3,3,3,3)
    19 = @(15,5,10,10,3,3)
    20 = @(12,5,10,10,3,3)
    21 = @(16,5,10,10,3,3)
    22 = @(12,5,10,10,3,3)
    23 = @(12,5,10,10,3,3)
    24 = @(12,5,10,10,3,3)
    25 = @(12,5,10,10,3,3)
    26 = @(12,5,10,10,3,3)
    27 = @(16,5,10,10,3,3)
    28 = @(12,5,10,10,3,3)
    29 = @(12,5,10,10,3,3)
    30 = @(14,5,10,10,3,3)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}

# Matches the author's final selection state (<selection activeCell="C4" sqref="C4"/>)
$ws.Range("C4").Select()
